# Commit: "#7 F-value in table 3 for Nodule model updated"
# Table 3 (Nodule count row, row 11) F-values are updated:
#   Treatment (B11):            11.33**  -> 10.87***
#   Group (C11):                --       -> 66.58***
#   Treatment x Group (D11):    --       -> 4.23**

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B11: Treatment F-value for Nodule count ---
$ws.Range("B11").Value = "10.87***"
$charsB = $ws.Range("B11").Characters(6, 3)
$charsB.Font.Superscript = $true

# --- C11: Group F-value for Nodule count (was a literal "--") ---
# Leading apostrophe keeps it quote-prefixed text (as the original "--" was).
$ws.Range("C11").Value = "'66.58***"
$charsC = $ws.Range("C11").Characters(6, 3)
$charsC.Font.Superscript = $true
# This cell now carries a distinct number format (matches the new style
# slot introduced for this cell in the saved workbook).
$ws.Range("C11").NumberFormat = "0.00"

# --- D11: Treatment x Group F-value for Nodule count (was a literal "--") ---
$ws.Range("D11").Value = "'4.23**"
$charsD = $ws.Range("D11").Characters(5, 2)
$charsD.Font.Superscript = $true

# Column C widened slightly to fit the new, longer values
# (target stored width 11.1640625; 10.33 is the closest input this host
# resolves to that stored width, since ColumnWidth here snaps to whole
# character counts before re-adding the fixed pixel padding).
$ws.Columns("C").ColumnWidth = 10.33

# Selection moved as part of the authoring session.
$ws.Range("E13").Select()
